# 10.4.1 worksheet update:
#  - add a 2023 data column (Q)
#  - add trilingual footnote text in row 6 (Kyrgyz / Russian / English) with a
#    superscript "1" marker, italic 8pt Times New Roman
#  - resize columns A:C to a uniform width, tweak a few row heights
#  - refresh the view zoom

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New "2023" column (Q): header + data point, copying the formatting that the
# neighbouring 2022 column (P) already uses.
# ---------------------------------------------------------------------------
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2023

$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 23.4

# ---------------------------------------------------------------------------
# Row 6 footnotes (Kyrgyz / Russian / English), each a two-run rich string:
# a superscript "1" followed by the italic footnote sentence. Alignment is
# applied before the per-character font tweaks so the whole cell (and its
# runs) end up sharing one tidy style.
# ---------------------------------------------------------------------------
$a6 = $ws.Range("A6")
$a6.Value = "1 2019-жылдан баштап маалыматтар, 2008 жылдагы Улуттук Эсептер Тутумунун эл аралык стандарттарына ылайык эсептелген "
$a6.HorizontalAlignment = -4131
$a6.VerticalAlignment = -4108
$a6.WrapText = $true
$a6r1 = $a6.Characters(1,1)
$a6r1.Font.Name = "Times New Roman"
$a6r1.Font.Size = 8
$a6r1.Font.Italic = $true
$a6r1.Font.Superscript = $true
$a6r1.Font.ColorIndex = 0
$a6r2 = $a6.Characters(2,115)
$a6r2.Font.Name = "Times New Roman"
$a6r2.Font.Size = 8
$a6r2.Font.Italic = $true
$a6r2.Font.ColorIndex = 0

$b6 = $ws.Range("B6")
$b6.Value = "1 Данные с 2019 года рассчитаны по международному стандарту Системы Национальных Счетов 2008 года"
$b6.HorizontalAlignment = -4131
$b6.VerticalAlignment = -4108
$b6.WrapText = $true
$b6r1 = $b6.Characters(1,2)
$b6r1.Font.Name = "Times New Roman"
$b6r1.Font.Size = 8
$b6r1.Font.Italic = $true
$b6r1.Font.Superscript = $true
$b6r1.Font.ColorIndex = 0
$b6r2 = $b6.Characters(3,95)
$b6r2.Font.Name = "Times New Roman"
$b6r2.Font.Size = 8
$b6r2.Font.Italic = $true
$b6r2.Font.ColorIndex = 0

$c6 = $ws.Range("C6")
$c6.Value = "1 Data from 2019 are calculated according to the international standard of the System of National Accounts 2008"
$c6.HorizontalAlignment = -4131
$c6.VerticalAlignment = -4108
$c6.WrapText = $true
$c6.NumberFormat = "@"
$c6r1 = $c6.Characters(1,1)
$c6r1.Font.Name = "Times New Roman"
$c6r1.Font.Size = 8
$c6r1.Font.Italic = $true
$c6r1.Font.Superscript = $true
$c6r1.Font.ColorIndex = 0
$c6r2 = $c6.Characters(2,110)
$c6r2.Font.Name = "Times New Roman"
$c6r2.Font.Size = 8
$c6r2.Font.Italic = $true
$c6r2.Font.ColorIndex = 0

$ws.Rows.Item(6).RowHeight = 36

# ---------------------------------------------------------------------------
# Column widths: A:C become a uniform, slightly narrower width.
# ---------------------------------------------------------------------------
$ws.Range("A:C").ColumnWidth = 36.71

# ---------------------------------------------------------------------------
# A few row-height tweaks that came along with the re-layout.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 53.25
$ws.Rows.Item(2).RowHeight = 14.25
$ws.Rows.Item(3).RowHeight = 14.25

# ---------------------------------------------------------------------------
# Refresh the view: normal zoom at 100%, selection back on A1.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 100
$ws.Range("A1").Select()
